# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (E2) and "Correspond Handback DateTime" (H2)
# values on the zh-cn and de-de sheets to reflect the latest handback run timestamps.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-17 10:29:45"
$wsZhCn.Range("H2").Value = "2016-03-17 10:30:00"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-17 10:29:49"
$wsDeDe.Range("H2").Value = "2016-03-17 10:30:11"
